$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.432.28"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "2.633.73"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'600.89"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "'153.05"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.561"
$ws.Range("E8").Value = "  +3.66%  "
$ws.Range("D9").Value = "2.634.75"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'0.123"
$ws.Range("E10").Value = "  +3.73%  "
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "'5.19"
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("D14").Value = "'27.69"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "3.112.52"
$ws.Range("E15").Value = "  +0.89%  "
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "67.474.41"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("D18").Value = "2.641.34"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'11.21"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "'365.32"
$ws.Range("E20").Value = "  +2.20%  "
$ws.Range("D21").Value = "'7.52"
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "'2.13"
$ws.Range("E23").Value = "  +4.39%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "'70.87"
$ws.Range("E25").Value = "  +4.86%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'10.27"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "2.763.90"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'0.0000103"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'580.20"
$ws.Range("E30").Value = "  -6.93%  "
$ws.Range("D31").Value = "'1.40"
$ws.Range("E31").Value = "  -3.67%  "
$ws.Range("D32").Value = "'7.83"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").Value = "'1.85"
$ws.Range("E33").Value = "  -0.20%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.128"
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("D37").Value = "'4.94"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("D38").Value = "'158.42"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("D39").Value = "'19.31"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'0.369"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Value = "'5.29"
$ws.Range("E41").Value = "  -3.20%  "
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  +3.44%  "
$ws.Range("D44").Value = "'41.24"
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "'16.37"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'156.01"
$ws.Range("E47").Value = "  +0.47%  "
$ws.Range("D48").Value = "0.0₆0286"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("E49").Value = "  -0.32%  "
$ws.Range("D50").Value = "'0.624"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'20.52"
$ws.Range("E51").Value = "  -1.83%  "
